# Bases update, 15-06-2024 21:10 - refreshes the trailing batch of matches
# (rows 300-306) in the "Poland Ekstraklasa" sheet with the latest
# teams / score / odds data. Column A (running id) is left untouched;
# every other column (B..AD) is rewritten per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row=300; B=7083187; E="Lech Poznan";            F="Korona Kielce";        G=1; H=2; I=1; J=0; K="A"; L=1.8;   M=3.8;  N=3.6;  O=2.1;  P=3.7;   Q=2.9;   R=-0.25; S=1.9;   T=1.95;  U=2.75; V=1.925; W=1.925; X=-1;                  Y=-1; Z=1.9;   AA=-1;                  AB=0.95;                AC=0.4625;              AD=-0.5 }
  @{ Row=301; B=7041338; E="Jagiellonia Bialystok";  F="Warta Poznan";         G=3; H=0; I=3; J=0; K="H"; L=1.444; M=4.75; N=5.25; O=1.4;  P=4.75;  Q=5.75;  R=-1.25; S=1.9;   T=1.95;  U=3;    V=1.925; W=1.925; X=0.3999999999999999; Y=-1; Z=-1;    AA=0.8999999999999999; AB=-1;                  AC=0;                   AD=0 }
  @{ Row=302; B=7088350; E="Puszcza Niepolomice";    F="Piast Gliwice";        G=1; H=0; I=0; J=0; K="H"; L=3;     M=3.1;  N=2.3;  O=2.7;  P=3;     Q=2.625; R=0;     S=1.975; T=1.875; U=2.25; V=2.025; W=1.825; X=1.7;                 Y=-1; Z=-1;    AA=0.9750000000000001; AB=-1;                  AC=-1;                  AD=0.825 }
  @{ Row=303; B=7083189; E="Pogon Szczecin";         F="Gornik Zabrze";        G=1; H=0; I=0; J=0; K="H"; L=1.727; M=4;    N=3.75; O=1.55; P=4.333; Q=4.5;   R=-1;    S=1.925; T=1.925; U=3.5;  V=2.025; W=1.825; X=0.55;                Y=-1; Z=-1;    AA=0;                   AB=0;                   AC=-1;                  AD=0.825 }
  @{ Row=304; B=7090293; E="Radomiak Radom";         F="Widzew Lodz";          G=1; H=3; I=1; J=0; K="A"; L=2.2;   M=3.1;  N=3.1;  O=2.15; P=3.2;   Q=3.1;   R=-0.25; S=1.925; T=1.925; U=2.75; V=1.9;   W=1.95;  X=-1;                  Y=-1; Z=2.1;   AA=-1;                  AB=0.925;               AC=0.8999999999999999; AD=-1 }
  @{ Row=305; B=7074364; E="Rakow Czestochowa";      F="Slask Wroclaw";        G=1; H=2; I=1; J=0; K="A"; L=2.5;   M=3.6;  N=2.4;  O=2.15; P=3.6;   Q=2.875; R=-0.25; S=1.95;  T=1.9;   U=2.5;  V=1.875; W=1.975; X=-1;                  Y=-1; Z=1.875; AA=-1;                  AB=0.8999999999999999; AC=0.875;               AD=-1 }
  @{ Row=306; B=7093821; E="LKS Lodz";                F="Stal Mielec";         G=3; H=2; I=3; J=0; K="H"; L=2.5;   M=3.4;  N=2.5;  O=2.2;  P=3.5;   Q=2.8;   R=-0.25; S=2.025; T=1.825; U=3;    V=2;     W=1.85;  X=1.2;                 Y=-1; Z=-1;    AA=1.025;               AB=-1;                  AC=1;                   AD=-1 }
)

foreach ($r in $rows) {
  $row = $r.Row
  $ws.Cells.Item($row, 2).Value  = $r.B   # B: match id
  $ws.Cells.Item($row, 5).Value  = $r.E   # E: HomeTeam
  $ws.Cells.Item($row, 6).Value  = $r.F   # F: AwayTeam
  $ws.Cells.Item($row, 7).Value  = $r.G   # G: FTHG
  $ws.Cells.Item($row, 8).Value  = $r.H   # H: FTAG
  $ws.Cells.Item($row, 9).Value  = $r.I   # I: HTHG
  $ws.Cells.Item($row, 10).Value = $r.J   # J: HTAG
  $ws.Cells.Item($row, 11).Value = $r.K   # K: FTR
  $ws.Cells.Item($row, 12).Value = $r.L   # L: oddH_op
  $ws.Cells.Item($row, 13).Value = $r.M   # M: oddD_op
  $ws.Cells.Item($row, 14).Value = $r.N   # N: oddA_op
  $ws.Cells.Item($row, 15).Value = $r.O   # O: oddH
  $ws.Cells.Item($row, 16).Value = $r.P   # P: oddD
  $ws.Cells.Item($row, 17).Value = $r.Q   # Q: oddA
  $ws.Cells.Item($row, 18).Value = $r.R   # R: Ah
  $ws.Cells.Item($row, 19).Value = $r.S   # S: oddAHH
  $ws.Cells.Item($row, 20).Value = $r.T   # T: oddAHA
  $ws.Cells.Item($row, 21).Value = $r.U   # U: AhOU
  $ws.Cells.Item($row, 22).Value = $r.V   # V: oddAHOver
  $ws.Cells.Item($row, 23).Value = $r.W   # W: oddAHUnder
  $ws.Cells.Item($row, 24).Value = $r.X   # X: PLH
  $ws.Cells.Item($row, 25).Value = $r.Y   # Y: PLD
  $ws.Cells.Item($row, 26).Value = $r.Z   # Z: PLA
  $ws.Cells.Item($row, 27).Value = $r.AA  # AA: PL_Ahh
  $ws.Cells.Item($row, 28).Value = $r.AB  # AB: PL_Aha
  $ws.Cells.Item($row, 29).Value = $r.AC  # AC: PL_AhOver
  $ws.Cells.Item($row, 30).Value = $r.AD  # AD: PL_AhUnder
}
